# Setup basic UI and added resources
# Clear the sample/placeholder row (row 5) that previously contained
# an example "Had some issues with..." description row, so the log
# starts with a single filled-in example row (row 4) followed by
# blank rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B5:D5").ClearContents()
